$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for "Ajo" / "Chino" / "Primera"
# at Mercado Mayorista Lo Valledor de Santiago. This pushes the existing
# data block (rows 770-818) down by one row (to 771-819) and inserts the
# new observation as the new row 770.

$ws.Rows.Item(770).Insert()

# Fill in the newly inserted row 770 with the static (repeated) columns
# plus the new weekly figures.
$ws.Range("A770").Value = 6
$ws.Range("B770").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C770").Value = "Metropolitana"
$ws.Range("D770").Value = 44826
$ws.Range("E770").Value = 13
$ws.Range("F770").Value = 100112003
$ws.Range("G770").Value = "Ajo"
$ws.Range("H770").Value = "Chino"
$ws.Range("I770").Value = "Primera"
$ws.Range("J770").Value = 2500
$ws.Range("K770").Value = 16000
$ws.Range("L770").Value = 17000
$ws.Range("M770").Value = 16600
$ws.Range("N770").Value = "`$/caja 10 kilos"
$ws.Range("O770").Value = "China"
$ws.Range("P770").Value = 1660
$ws.Range("Q770").Value = 10
$ws.Range("R770").Value = "Hortaliza"
